$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - C1 now references "runMode" (shared string index changes
# automatically because unused strings are dropped on save)
# Row 1 values are unchanged (username / password / runMode) - nothing to do.

# Update row 2 data values (new sample credentials)
$ws.Range("A2").Value = "mngr289535"
$ws.Range("B2").Value = "AtUpypU"
$ws.Range("C2").Value = "y"

# Remove the hyperlinks on B4 and B5 (keeps the Hyperlink cell style)
$ws.Range("B4").Hyperlinks.Delete()
$ws.Range("B5").Hyperlinks.Delete()

# Row 3 is removed entirely (no longer used)
$ws.Range("A3:C3").Clear()

# Rows 4 and 5 only keep the (now empty) B cell with its hyperlink style
$ws.Range("A4").Clear()
$ws.Range("C4").Clear()
$ws.Range("B4").ClearContents()

$ws.Range("A5").Clear()
$ws.Range("C5").Clear()
$ws.Range("B5").ClearContents()

# Update the active selection
$ws.Range("B11").Select()

$wb.Save()
